$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.552.99"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.880.68"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.23"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4718"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2884"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06530"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.02"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "101.10"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07818"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7428"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "1.876.27"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.206"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "285.38"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "30.523.94"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.14"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9993"
$ws.Range("D21").Value = "2.119.61"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.363"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9990"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.379"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.106"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.51"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09696"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.323"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.496"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.254"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.181"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04834"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6932"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.772"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01904"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.862"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.24"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.337"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.975"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4242"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9987"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8284"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.33"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.783"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.043"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05765"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "889.04"
$ws.Range("E51").Value = "  -3.32%  "
